# Fixing Tour create template
#
# - The "Note" label moves from B16 up to its own row, B15.
# - The "Ticket Description" label moves from I16 to H17, leaving I16 blank
#   (but still formatted).
# - A new "Route available" label is added at B17, alongside the relocated
#   "Ticket Description" at H17.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Propagate the bold label formatting (currently on B16 / I16) to the new
# homes for these labels before we touch their contents.
$ws.Range("B16").Copy()
$ws.Range("B15").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("I16").Copy()
$ws.Range("B17").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("I16").Copy()
$ws.Range("H17").PasteSpecial(-4122)   # xlPasteFormats

$excel.CutCopyMode = $false

# Move the note text from B16 up to B15, then remove B16 entirely.
$ws.Range("B15").Value = "Note: Scanning from number 1 -> 10. Please fill out the whole row to create a tour"
$ws.Range("B16").Clear()

# Relocate "Ticket Description" from I16 to H17, leaving I16 blank (but
# still carrying its original formatting).
$ws.Range("H17").Value = "Ticket Description"
$ws.Range("I16").Value = $null

# Add the new "Route available" label.
$ws.Range("B17").Value = "Route available"

# Update the active selection to match the saved workbook state.
$ws.Range("F19").Select()
